# Generate Report for Handoff
# The file "6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e.md" has been handed off again:
#   - Overview sheet row for this file: zh-cn & de-de status -> "Ready for handoff",
#     Latest HO Xliff Generate Date -> 2016-08-22 06:47:18
#   - zh-cn sheet row for this file: Status -> "Ready for handoff",
#     Latest Handoff Datetime -> 2016-08-22 06:47:14,
#     Error Detail -> stale-handback-file warning message
#   - de-de sheet row for this file: Status -> "Ready for handoff",
#     Latest Handoff Datetime -> 2016-08-22 06:47:18,
#     Error Detail -> stale-handback-file warning message
#   - zh-cn / de-de "Error Detail" column widened to fit the new message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39456a458bb60bf5da5a9f398411b755ed8795ca/e2e/6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02a06d3a42f4cb81c8b14836f09e0d4a194fdf20/e2e/6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e.md."

# --- Overview sheet: row 3 is the 6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-22 06:47:18"

# --- zh-cn sheet: row 3 is the 6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e.md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-22 06:47:14"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is the 6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e.md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-22 06:47:18"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664

Write-Output "Report regenerated for handoff."
